$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices with dotted thousand separators,
# trailing zeros, scientific-looking small decimals) are stored as literal text,
# matching the source data which treats these as inline strings, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.250.81"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.135.52"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.71"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.12"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +27.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.374"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.131.13"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.772"
$ws.Range("E11").Value = "  +26.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.197"
$ws.Range("E12").Value = "  +8.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.89"
$ws.Range("E14").Value = "  +8.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.59"
$ws.Range("E15").Value = "  +5.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.083.49"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.708.89"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.121.62"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.77"
$ws.Range("E19").Value = "  +11.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000216"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.20"
$ws.Range("E21").Value = "  +5.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.81"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.88"
$ws.Range("E23").Value = "  +7.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.23"
$ws.Range("E24").Value = "  +5.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.34"
$ws.Range("E25").Value = "  +14.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.55"
$ws.Range("E26").Value = "  +6.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.38"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.293.82"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.165"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.23"
$ws.Range("E31").Value = "  +13.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.33"
$ws.Range("E32").Value = "  +4.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.907"
$ws.Range("E33").Value = "  -14.38%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.22"
$ws.Range("E34").Value = "  +12.54%  "
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.72"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.146"
$ws.Range("E36").Value = "  +14.16%  "
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("E38").Value = "  +5.34%  "
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.163"
$ws.Range("E40").Value = "  +20.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.27"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0838"
$ws.Range("E43").Value = "  +19.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.406"
$ws.Range("E44").Value = "  +11.60%  "
$ws.Range("E45").Value = "  +6.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.37"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.19"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.32"
$ws.Range("E49").Value = "  +8.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.40"
$ws.Range("E50").Value = "  +10.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "171.49"
$ws.Range("E51").Value = "  +7.53%  "
